$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "44.088.39"
$ws.Range("E2").Value = "  +3.02%  "

# Row 3
$ws.Range("D3").Value = "2.278.30"
$ws.Range("E3").Value = "  +3.07%  "

# Row 4
$ws.Range("E4").Value = "  -0.25%  "

# Row 5
$ws.Range("D5").Value = "'318.67"
$ws.Range("E5").Value = "  +1.43%  "

# Row 6
$ws.Range("E6").Value = "  +8.18%  "

# Row 7
$ws.Range("E7").Value = "  +1.52%  "

# Row 8
$ws.Range("E8").Value = "  -0.20%  "

# Row 9
$ws.Range("D9").Value = "'0.574"
$ws.Range("E9").Value = "  +3.00%  "

# Row 10
$ws.Range("D10").Value = "'39.03"
$ws.Range("E10").Value = "  +7.37%  "

# Row 11
$ws.Range("E11").Value = "  +2.15%  "

# Row 12
$ws.Range("D12").Value = "'7.90"
$ws.Range("E12").Value = "  +2.15%  "

# Row 13
$ws.Range("E13").Value = "  +1.95%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.884"

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.624.75"
$ws.Range("E15").Value = "  +3.11%  "

# Row 16
$ws.Range("D16").Value = "'14.67"
$ws.Range("E16").Value = "  +3.97%  "

# Row 17
$ws.Range("D17").Value = "2.280.71"
$ws.Range("E17").Value = "  +3.59%  "

# Row 18
$ws.Range("D18").Value = "43.996.32"
$ws.Range("E18").Value = "  +3.04%  "

# Row 19
$ws.Range("D19").Value = "'14.05"
$ws.Range("E19").Value = "  -6.08%  "

# Row 20
$ws.Range("D20").Value = "0.0₂01000"
$ws.Range("E20").Value = "  +4.47%  "

# Row 21
$ws.Range("E21").Value = "  +3.03%  "

# Row 22
$ws.Range("D22").Value = "'66.25"
$ws.Range("E22").Value = "  +1.96%  "

# Row 23
$ws.Range("E23").Value = "  +1.97%  "

# Row 24
$ws.Range("D24").Value = "'238.13"
$ws.Range("E24").Value = "  +1.35%  "

# Row 25
$ws.Range("E25").Value = "  +4.70%  "

# Row 26
$ws.Range("E26").Value = "  -0.09%  "

# Row 27
$ws.Range("D27").Value = "'10.28"
$ws.Range("E27").Value = "  +2.76%  "

# Row 28
$ws.Range("D28").Value = "'38.87"
$ws.Range("E28").Value = "  +15.15%  "

# Row 29
$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = "  -0.83%  "

# Row 30
$ws.Range("D30").Value = "'6.52"
$ws.Range("E30").Value = "  +4.43%  "

# Row 31
$ws.Range("D31").Value = "'163.58"
$ws.Range("E31").Value = "  +4.58%  "

# Row 32
$ws.Range("D32").Value = "'20.63"
$ws.Range("E32").Value = "  +1.33%  "

# Row 33
$ws.Range("D33").Value = "'0.0888"
$ws.Range("E33").Value = "  +2.28%  "

# Row 34
$ws.Range("E34").Value = "  -1.09%  "

# Row 35
$ws.Range("D35").Value = "'2.09"
$ws.Range("E35").Value = "  +5.60%  "

# Row 36
$ws.Range("D36").Value = "'3.26"
$ws.Range("E36").Value = "  +2.39%  "

# Row 37
$ws.Range("D37").Value = "'0.114"
$ws.Range("E37").Value = "  +11.66%  "

# Row 38
$ws.Range("E38").Value = "  -0.06%  "

# Row 39
$ws.Range("D39").Value = "'4.00"
$ws.Range("E39").Value = "  +8.96%  "

# Row 40
$ws.Range("E40").Value = "  +1.93%  "

# Row 41
$ws.Range("D41").Value = "'15.50"
$ws.Range("E41").Value = "  +27.97%  "

# Row 42
$ws.Range("D42").Value = "'0.0327"
$ws.Range("E42").Value = "  +1.06%  "

# Row 43
$ws.Range("E43").Value = "  -0.08%  "

# Row 44
$ws.Range("D44").Value = "1.767.37"
$ws.Range("E44").Value = "  -6.11%  "

# Row 45
$ws.Range("E45").Value = "  +1.35%  "

# Row 46
$ws.Range("D46").Value = "'85.79"
$ws.Range("E46").Value = "  -3.02%  "

# Row 47
$ws.Range("D47").Value = "'5.42"
$ws.Range("E47").Value = "  +0.72%  "

# Row 48
$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").Value = "'59.89"
$ws.Range("E48").Value = "  -0.48%  "

# Row 49
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").Value = "'75.30"
$ws.Range("E49").Value = "  +0.77%  "

# Row 50
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'8.79"
$ws.Range("E50").Value = "  +3.24%  "

# Row 51
$ws.Range("D51").Value = "'104.55"
$ws.Range("E51").Value = "  +3.61%  "
